$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# --- Move the hyperlink from G3 to G2 (same target mailto address, new location) ---
# NB: Range.Hyperlinks.Delete() clears the whole sheet's hyperlink collection on this
# host, so it must run before the new one is added (there is only one link to begin
# with, so this is safe).
$ws.Range("G3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:aashish.kumar@sofbang.com;%20false")

# Re-apply G3's original Hyperlink-style formatting onto G2 (reuse the existing style
# rather than letting Hyperlinks.Add mint a new near-duplicate one) - do this last so
# it isn't clobbered by the Add call above.
$ws.Range("G3").Copy()
$ws.Range("G2").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 2 (Get User Account): fill in ResponseMapKeys / ResponseMapValues ---
$ws.Range("F2").Value = "account.email_address;account.is_locked;account.quotas.api_signature_requests_left"
$ws.Range("G2").Value = "aashish.kumar@sofbang.com;false;5000"

# --- Row 3 (Verify User Account): clear the old ResponseMapKeys / ResponseMapValues ---
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""

# --- Row 5: rename "Verify Team" -> "Get Team" and populate its ResponseMapKeys / ResponseMapValues ---
$ws.Range("B5").Value = "Get Team"
$ws.Range("F5").Value = "team.name;team.accounts[0].email_address"
$ws.Range("G5").Value = "Sofbang Team;vivek.ahuja@sofbang.com"

# --- Selection moves to G6 in the saved view ---
$ws.Range("G6").Select()
